{"js": "const body = context.document.body;\n\n// Delete the lone-space run between \")\" and \" 1 ni\" by merging it into the \")\" run.\n{\n  const r = body.search(\") \", { matchCase: true });\n  r.load(\"text\");\n  await context.sync();\n  r.items[0].insertText(\")\", Word.InsertLocation.replace);\n}\n\n// Text replacements (Swahili -> English), one per originally-unique run.\nconst replacements = [\n  [\"Wafungwa na peremende - manukuu:\", \"Prisoners and candies - subtitles:\"],\n  [\"**mazungumzo huanza saa 55 ya pili sio 27 kwa sababu ya klipu ya utangulizi. Nilirekebisha nyakati ipasavyo. -John Argentino\", \"**dialogue starts at second 55 not 27 because of the intro clip. I adjusted the times accordingly. -John Argentino\"],\n  [\"[Muziki]\", \"[Music]\"],\n  [\"wanahisabati wanne mkali wanachukuliwa\", \"four bright mathematicians are taken into\"],\n  [\"chini ya ulinzi na kuwekwa gerezani kwa sababu walijaribu\", \"custody and put in jail because they tried\"],\n  [\"kumshawishi mwanamke mzee kuwa Goedel's\", \"to convince an old lady that the Goedel's\"],\n  [\"nadharia za kutokamilika ni kweli. Kila\", \"incompleteness theorems are true. Every\"],\n  [\"mtaalamu wa hisabati ana kiini chake ambacho sisi\", \"mathematician has his own cell that we\"],\n  [\"inaweza kuhesabu na nambari kutoka 1 hadi 4.\", \"can enumerate with a number from 1 to 4.\"],\n  [\"kabla ya kuingia kwenye seli fulani\", \"before entering the cell a certain\"],\n  [\"idadi ya peremende kubwa kuliko \", \"number of candies greater than \"],\n  [\"e: AU \", \"e: OR \"],\n  [\"SAWA NA\", \"EQUAL TO\"],\n  [\" 1 ni\", \" 1 is\"],\n  [\"wanapewa kila mtaalamu wa hisabati na wao\", \"given to every mathematician and they\"],\n  [\"wanaambiwa wana peremende 11 kwa jumla.\", \"are told they have 11 candies in total.\"],\n  [\"lakini kila mtu anajua idadi yake tu\", \"but everyone knows only his number of\"],\n  [\"pipi na jumla. 1 na sio\", \"candies and the total. 1 and is not\"],\n  [\"kuruhusiwa kuuliza nambari zingine.\", \"allowed to ask for the others number.\"],\n  [\"kisha mwanahisabati wa kwanza anauliza\", \"then the first mathematician asks the\"],\n  [\"pili: 'namba 2 unajua kama wewe\", \"second: 'number 2 do you know if you\"],\n  [\"kuwa na peremende nyingi kuliko mimi?' ya pili\", \"have more candies than me?' the second\"],\n  [\"mwanahisabati anajibu hana. Kisha yeye\", \"mathematician answers he doesn't. Then he\"],\n  [\"anauliza kwa nambari 3: 'unajua kama unayo\", \"asks to number 3: 'do you know if you have\"],\n  [\"pipi zaidi kuliko mimi?'\", \"more candy than me?'\"],\n  [\"mwanahisabati wa tatu anajibu: 'hapana niko\", \"the third mathematician answers: 'no I'm\"],\n  [\"samahani sifanyi'. Katika hatua hii ya nne\", \"sorry I don't'. At this point the fourth\"],\n  [\"mtaalamu wa hisabati anasema: 'jamani mnafahamu\", \"mathematician says: 'hey guys you know\"],\n  [\"nini, najua hasa pipi ngapi\", \"what, I know exactly how many candies\"],\n  [\"kila mtu ana hapa'. Cha kushangaza hata\", \"everyone has here'. Surprisingly even the\"],\n  [\"wanahisabati wengine watatu wanasema hivyo sasa\", \"other three mathematicians say that now\"],\n  [\"wanajua kila mtu ana pipi ngapi\", \"they know how many candies everyone has\"],\n  [\"kwa hivyo swali ni: unaweza kujua\", \"so the question is: can you figure out\"],\n  [\"idadi ya pipi kila mfungwa ana\", \"the number of candies every prisoner has\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nreturn \"ok\";", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-AllOccurrences($oldText, $newText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.Text = $oldText\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop - never wrap back to the start\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $guard = 0\n    while ($find.Execute()) {\n        $rng.Text = $newText\n        $guard++\n        if ($guard -gt 50) { break }\n    }\n}\n\n# Delete the lone-space run between \")\" and \" 1 ni\" by merging it into the \")\" run.\nReplace-AllOccurrences ') ' ')'\n\n# Text replacements (Swahili -> English), one per originally-unique run.\n$replacements = @(\n    @('Wafungwa na peremende - manukuu:', 'Prisoners and candies - subtitles:'),\n    @('**mazungumzo huanza saa 55 ya pili sio 27 kwa sababu ya klipu ya utangulizi. Nilirekebisha nyakati ipasavyo. -John Argentino', '**dialogue starts at second 55 not 27 because of the intro clip. I adjusted the times accordingly. -John Argentino'),\n    @('[Muziki]', '[Music]'),\n    @('wanahisabati wanne mkali wanachukuliwa', 'four bright mathematicians are taken into'),\n    @('chini ya ulinzi na kuwekwa gerezani kwa sababu walijaribu', 'custody and put in jail because they tried'),\n    @('kumshawishi mwanamke mzee kuwa Goedel''s', 'to convince an old lady that the Goedel''s'),\n    @('nadharia za kutokamilika ni kweli. Kila', 'incompleteness theorems are true. Every'),\n    @('mtaalamu wa hisabati ana kiini chake ambacho sisi', 'mathematician has his own cell that we'),\n    @('inaweza kuhesabu na nambari kutoka 1 hadi 4.', 'can enumerate with a number from 1 to 4.'),\n    @('kabla ya kuingia kwenye seli fulani', 'before entering the cell a certain'),\n    @('idadi ya peremende kubwa kuliko ', 'number of candies greater than '),\n    @('e: AU ', 'e: OR '),\n    @('SAWA NA', 'EQUAL TO'),\n    @(' 1 ni', ' 1 is'),\n    @('wanapewa kila mtaalamu wa hisabati na wao', 'given to every mathematician and they'),\n    @('wanaambiwa wana peremende 11 kwa jumla.', 'are told they have 11 candies in total.'),\n    @('lakini kila mtu anajua idadi yake tu', 'but everyone knows only his number of'),\n    @('pipi na jumla. 1 na sio', 'candies and the total. 1 and is not'),\n    @('kuruhusiwa kuuliza nambari zingine.', 'allowed to ask for the others number.'),\n    @('kisha mwanahisabati wa kwanza anauliza', 'then the first mathematician asks the'),\n    @('pili: ''namba 2 unajua kama wewe', 'second: ''number 2 do you know if you'),\n    @('kuwa na peremende nyingi kuliko mimi?'' ya pili', 'have more candies than me?'' the second'),\n    @('mwanahisabati anajibu hana. Kisha yeye', 'mathematician answers he doesn''t. Then he'),\n    @('anauliza kwa nambari 3: ''unajua kama unayo', 'asks to number 3: ''do you know if you have'),\n    @('pipi zaidi kuliko mimi?''', 'more candy than me?'''),\n    @('mwanahisabati wa tatu anajibu: ''hapana niko', 'the third mathematician answers: ''no I''m'),\n    @('samahani sifanyi''. Katika hatua hii ya nne', 'sorry I don''t''. At this point the fourth'),\n    @('mtaalamu wa hisabati anasema: ''jamani mnafahamu', 'mathematician says: ''hey guys you know'),\n    @('nini, najua hasa pipi ngapi', 'what, I know exactly how many candies'),\n    @('kila mtu ana hapa''. Cha kushangaza hata', 'everyone has here''. Surprisingly even the'),\n    @('wanahisabati wengine watatu wanasema hivyo sasa', 'other three mathematicians say that now'),\n    @('wanajua kila mtu ana pipi ngapi', 'they know how many candies everyone has'),\n    @('kwa hivyo swali ni: unaweza kujua', 'so the question is: can you figure out'),\n    @('idadi ya pipi kila mfungwa ana', 'the number of candies every prisoner has'),\n)\n\nforeach ($pair in $replacements) {\n    Replace-AllOccurrences $pair[0] $pair[1]\n}"}
